$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ACHENGLI LAILA"
$ws.Range("B2").Value = "J207703"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "00101211115087750001201090"
$ws.Range("D2").Value = "Ait souss"
$ws.Range("E2").Value = "BP Centre Sud"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "901/FES "
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 5000
$ws.Range("J2").Value = 500
$ws.Range("K2").Value = 4500

# Row 3
$ws.Range("A3").Value = "CHARIJI ABDELLAH"
$ws.Range("B3").Value = "BJ36877"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "00101211111292695000201732"
$ws.Range("D3").Value = "AOURIR"
$ws.Range("E3").Value = "BP CENTRE SUD"
$ws.Range("F3").Value = "Logement de fonction"
$ws.Range("G3").Value = "901/LF/FES "
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 6000
$ws.Range("J3").Value = 600
$ws.Range("K3").Value = 5400

# Row 4 (totals row) - only numeric columns change
$ws.Range("I4").Value = 11000
$ws.Range("J4").Value = 1100
$ws.Range("K4").Value = 9900

$wb.Save()
